$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (2014/12) ---
$ws.Range("D2").Value = 5538
$ws.Range("E2").Value = 194
$ws.Range("F2").Value = 194
$ws.Range("G2").Value = 169
$ws.Range("H2").Value = 150
$ws.Range("I2").Value = 150
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value = 5787
$ws.Range("L2").Value = 1813
$ws.Range("M2").Value = 3974
$ws.Range("N2").Value = 3974
$ws.Range("O2").ClearContents()
$ws.Range("P2").Value = 346
$ws.Range("Q2").Value = 106
$ws.Range("R2").Value = -170
$ws.Range("S2").Value = 155
$ws.Range("T2").Value = 171
$ws.Range("U2").Value = -66
$ws.Range("V2").Value = 1263
$ws.Range("W2").Value = 3.51
$ws.Range("X2").Value = 2.72
$ws.Range("Y2").Value = 3.74
$ws.Range("Z2").Value = 2.64
$ws.Range("AA2").Value = 45.62
$ws.Range("AB2").Value = 921.4299999999999
$ws.Range("AC2").Value = 2174
$ws.Range("AD2").Value = 8.19
$ws.Range("AE2").Value = 57438
$ws.Range("AF2").Value = 0.31
$ws.Range("AG2").Value = 700
$ws.Range("AH2").Value = 3.93
$ws.Range("AI2").Value = 32.2
$ws.Range("AJ2").Value = 6918617

# --- Row 3 (2015/12) ---
$ws.Range("D3").Value = 4385
$ws.Range("E3").Value = 82
$ws.Range("F3").Value = 82
$ws.Range("G3").Value = 54
$ws.Range("H3").Value = 7
$ws.Range("I3").Value = 7
$ws.Range("J3").ClearContents()
$ws.Range("K3").Value = 5140
$ws.Range("L3").Value = 1203
$ws.Range("M3").Value = 3937
$ws.Range("N3").Value = 3937
$ws.Range("O3").ClearContents()
$ws.Range("P3").Value = 346
$ws.Range("Q3").Value = 1019
$ws.Range("R3").Value = -350
$ws.Range("S3").Value = -692
$ws.Range("T3").Value = 332
$ws.Range("U3").Value = 686
$ws.Range("V3").Value = 640
$ws.Range("W3").Value = 1.88
$ws.Range("X3").Value = 0.16
$ws.Range("Y3").Value = 0.18
$ws.Range("Z3").Value = 0.13
$ws.Range("AA3").Value = 30.55
$ws.Range("AB3").Value = 913.77
$ws.Range("AC3").Value = 100
$ws.Range("AD3").Value = 153.7
$ws.Range("AE3").Value = 56902
$ws.Range("AF3").Value = 0.27
$ws.Range("AG3").Value = 400
$ws.Range("AH3").Value = 2.6
$ws.Range("AI3").Value = 399.22
$ws.Range("AJ3").Value = 6918617

# --- Row 4 (2016/12) ---
$ws.Range("D4").Value = 3640
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = 13
$ws.Range("G4").Value = -7
$ws.Range("H4").Value = -6
$ws.Range("I4").Value = -6
$ws.Range("J4").ClearContents()
$ws.Range("K4").Value = 5868
$ws.Range("L4").Value = 2009
$ws.Range("M4").Value = 3859
$ws.Range("N4").Value = 3859
$ws.Range("O4").ClearContents()
$ws.Range("P4").Value = 346
$ws.Range("Q4").Value = -230
$ws.Range("R4").Value = -96
$ws.Range("S4").Value = 638
$ws.Range("T4").Value = 42
$ws.Range("U4").Value = -272
$ws.Range("V4").Value = 1249
$ws.Range("W4").Value = 0.37
$ws.Range("X4").Value = -0.15
$ws.Range("Y4").Value = -0.14
$ws.Range("Z4").Value = -0.1
$ws.Range("AA4").Value = 52.07
$ws.Range("AB4").Value = 902.22
$ws.Range("AC4").Value = -80
$ws.Range("AD4").Value = -191.55
$ws.Range("AE4").Value = 55779
$ws.Range("AF4").Value = 0.28
$ws.Range("AG4").Value = 500
$ws.Range("AH4").Value = 3.25
$ws.Range("AI4").Value = -621.9
$ws.Range("AJ4").Value = 6918617

# --- Row 5 (2017/12) ---
$ws.Range("D5").Value = 6906
$ws.Range("E5").Value = 227
$ws.Range("F5").Value = 227
$ws.Range("G5").Value = 87
$ws.Range("H5").Value = 37
$ws.Range("I5").Value = 37
$ws.Range("J5").ClearContents()
$ws.Range("K5").Value = 6591
$ws.Range("L5").Value = 2822
$ws.Range("M5").Value = 3769
$ws.Range("N5").Value = 3769
$ws.Range("O5").ClearContents()
$ws.Range("P5").Value = 346
$ws.Range("Q5").Value = -926
$ws.Range("R5").Value = 11
$ws.Range("S5").Value = 706
$ws.Range("T5").Value = 63
$ws.Range("U5").Value = -989
$ws.Range("V5").Value = 2004
$ws.Range("W5").Value = 3.29
$ws.Range("X5").Value = 0.53
$ws.Range("Y5").Value = 0.97
$ws.Range("Z5").Value = 0.59
$ws.Range("AA5").Value = 74.86
$ws.Range("AB5").Value = 906.73
$ws.Range("AC5").Value = 533
$ws.Range("AD5").Value = 29
$ws.Range("AE5").Value = 54480
$ws.Range("AF5").Value = 0.28
$ws.Range("AG5").Value = 500
$ws.Range("AH5").Value = 3.24
$ws.Range("AI5").Value = 93.87
$ws.Range("AJ5").Value = 6918617

# --- Row 6 (2018/12) ---
$ws.Range("D6").Value = 5888
$ws.Range("E6").Value = 179
$ws.Range("F6").Value = 179
$ws.Range("G6").Value = 197
$ws.Range("H6").Value = 172
$ws.Range("I6").Value = 172
$ws.Range("K6").Value = 6035
$ws.Range("L6").Value = 1974
$ws.Range("M6").Value = 4060
$ws.Range("N6").Value = 4060
$ws.Range("P6").Value = 377
$ws.Range("Q6").Value = 1139
$ws.Range("R6").Value = 71
$ws.Range("S6").Value = -909
$ws.Range("T6").Value = 85
$ws.Range("U6").Value = 1054
$ws.Range("V6").Value = 1138
$ws.Range("W6").Value = 3.05
$ws.Range("X6").Value = 2.92
$ws.Range("Y6").Value = 4.39
$ws.Range("Z6").Value = 2.72
$ws.Range("AA6").Value = 48.63
$ws.Range("AB6").Value = 888.5599999999999
$ws.Range("AC6").Value = 2362
$ws.Range("AD6").Value = 4.89
$ws.Range("AE6").Value = 53885
$ws.Range("AF6").Value = 0.21
$ws.Range("AG6").Value = 400
$ws.Range("AH6").Value = 3.46
$ws.Range("AI6").Value = 17.54
$ws.Range("AJ6").Value = 7534848

# --- Rows 7, 8, 9 (2019/12(E), 2020/12(E), 2021/12(E)) ---
# All forecast data cells (D..AI) are removed, keeping only A/B/C identifier cells.
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()
